$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: "_old" suffix columns become "_FV2404", "_new" suffix columns become "_FV2410"
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = ($baseNames[$i] + "_FV2404")
}

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = ($baseNames[$i] + "_FV2410")
}

# Turn the used range into an Excel Table (ListObject)
$rng = $ws.Range("A1:U90")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (top row)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
